$wb = $excel.ActiveWorkbook

$wsMuscles = $wb.Worksheets.Item("Muscles(L)")
$wsTendons = $wb.Worksheets.Item("Tendons(L)")

# Row 2
$wsMuscles.Range("B2").Value = 146
$wsMuscles.Range("C2").Value = 140
$wsMuscles.Range("D2").Value = 129

# Row 8
$wsMuscles.Range("B8").Value = 255
$wsMuscles.Range("C8").Value = 260
$wsMuscles.Range("D8").Value = 263
$wsMuscles.Range("E8").Value = 258
$wsMuscles.Range("F8").Value = 257
$wsMuscles.Range("G8").Value = 262

# Row 11
$wsMuscles.Range("B11").Value = 264
$wsMuscles.Range("C11").Value = 255
$wsMuscles.Range("D11").Value = 248

# Row 14
$wsMuscles.Range("C14").Value = 206

# Row 17
$wsMuscles.Range("C17").Value = 153
$wsMuscles.Range("D17").Value = 152
$wsMuscles.Range("E17").Value = 175

# Row 26
$wsMuscles.Range("B26").Value = 214
$wsMuscles.Range("C26").Value = 220
$wsMuscles.Range("D26").Value = 214
$wsMuscles.Range("E26").Value = 220

# Row 29
$wsMuscles.Range("C29").Value = 227
$wsMuscles.Range("D29").Value = 227

# Update the saved selection on Muscles(L) without leaving it as the active tab
$wsMuscles.Range("J28").Select()

# Tendons(L) sheet: fix duplicate shared string reference for "d"
$wsTendons.Range("I25").Value = "d"

# Restore Tendons(L) as the active/selected sheet (matches original tabSelected state)
$wsTendons.Select()
